$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 752.5714
$ws.Range("J19").Value = 748.8461
$ws.Range("L19").Value = 748.8461
$ws.Range("N19").Value = -1098.8461
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
# Row 38
$ws.Range("H38").Value = 548
$ws.Range("I38").Value = 448.25
$ws.Range("K38").Value = 1344.75
$ws.Range("M38").Value = -972.75
# Row 107
$ws.Range("H107").Value = 344.9
$ws.Range("I107").Value = 327.66666
$ws.Range("K107").Value = 327.66666
$ws.Range("M107").Value = 1592.33334
# Row 109
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774
# Row 112
$ws.Range("H112").Value = 1748.3334
$ws.Range("J112").Value = 2125
$ws.Range("L112").Value = 6375
$ws.Range("N112").Value = -8591
# Row 137
$ws.Range("H137").Value = 1500
$ws.Range("I137").Value = 1500
$ws.Range("K137").Value = 4500
$ws.Range("M137").Value = -1950

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 33
$ws.Range("H33").Value = 2000
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1671
$ws.Range("N33").ClearContents()
# Row 61
$ws.Range("H61").Value = 1659.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 1659.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 1659.5
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2083.5
# Row 74
$ws.Range("H74").Value = 1204.5
$ws.Range("I74").Value = 1204.5
$ws.Range("K74").Value = 1204.5
$ws.Range("M74").Value = -330.5
# Row 77
$ws.Range("H77").Value = 1204.5
$ws.Range("I77").Value = 1204.5
$ws.Range("K77").Value = 6022.5
$ws.Range("M77").Value = -1654.5
# Row 136
$ws.Range("H136").Value = 1659.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1659.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 4978.5
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -10078.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2384.1538
$ws.Range("I94").Value = 1999.1428
$ws.Range("J94").Value = 2833.3333
$ws.Range("K94").Value = 1999.1428
$ws.Range("L94").Value = 2833.3333
$ws.Range("M94").Value = -1548.1428
$ws.Range("N94").Value = -3735.3333
# Row 134
$ws.Range("H134").Value = 2406
$ws.Range("I134").Value = 812
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 2436
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = 99
$ws.Range("N134").Value = -17070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 44
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35884
# Row 58
$ws.Range("H58").Value = 1745
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 62
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376
# Row 65
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880
# Row 93
$ws.Range("H93").Value = 24740.666
$ws.Range("I93").Value = 24740.666
$ws.Range("K93").Value = 24740.666
$ws.Range("M93").Value = -22868.666
# Row 103
$ws.Range("H103").Value = 49000
$ws.Range("I103").Value = 49000
$ws.Range("K103").Value = 49000
$ws.Range("M103").Value = -47828
# Row 106
$ws.Range("H106").Value = 108333
$ws.Range("J106").Value = 108333
$ws.Range("L106").Value = 108333
$ws.Range("N106").Value = -110857
# Row 107
$ws.Range("H107").Value = 422.66666
$ws.Range("I107").Value = 422.66666
$ws.Range("K107").Value = 422.66666
$ws.Range("M107").Value = 1497.33334
# Row 122
$ws.Range("H122").Value = 1502.5
$ws.Range("I122").Value = 1670.3334
$ws.Range("J122").Value = 999
$ws.Range("K122").Value = 5011.0002
$ws.Range("L122").Value = 2997
$ws.Range("M122").Value = -2561.0002
$ws.Range("N122").Value = -7897
# Row 136
$ws.Range("H136").Value = 1745
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 12.034483
$ws.Range("I2").Value = 7.047619
$ws.Range("K2").Value = 42.285714
$ws.Range("M2").Value = 70.714286
# Row 18
$ws.Range("H18").Value = 3715
$ws.Range("I18").Value = 3715
$ws.Range("K18").Value = 11145
$ws.Range("M18").Value = -10976
# Row 23
$ws.Range("H23").Value = 110.666664
$ws.Range("I23").Value = 90
$ws.Range("J23").Value = 114.8
$ws.Range("K23").Value = 270
$ws.Range("L23").Value = 344.4
$ws.Range("M23").Value = -35
$ws.Range("N23").Value = -814.4
# Row 38
$ws.Range("H38").Value = 851.4
$ws.Range("I38").Value = 137.5
$ws.Range("J38").Value = 1029.875
$ws.Range("K38").Value = 412.5
$ws.Range("L38").Value = 3089.625
$ws.Range("M38").Value = -65.5
$ws.Range("N38").Value = -3783.625
# Row 109
$ws.Range("H109").Value = 813.5
$ws.Range("J109").Value = 900
$ws.Range("L109").Value = 2700
$ws.Range("N109").Value = -4780
# Row 131
$ws.Range("H131").Value = 4508
$ws.Range("J131").Value = 4361
$ws.Range("L131").Value = 13083
$ws.Range("N131").Value = -23163
# Row 133
$ws.Range("H133").Value = 400
$ws.Range("I133").Value = 400
$ws.Range("K133").Value = 1200
$ws.Range("M133").Value = 3860

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 1000
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -2996
# Row 83
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 5000
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -14984
# Row 102
$ws.Range("H102").Value = 58991.715
$ws.Range("I102").Value = 81938.60000000001
$ws.Range("J102").Value = 1624.5
$ws.Range("K102").Value = 81938.60000000001
$ws.Range("L102").Value = 1624.5
$ws.Range("M102").Value = -80316.60000000001
$ws.Range("N102").Value = -4868.5
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 126
$ws.Range("H126").Value = 500
$ws.Range("I126").Value = 500
$ws.Range("K126").Value = 1500
$ws.Range("M126").Value = 970

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1377.2
$ws.Range("I7").Value = 1221.5
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 1221.5
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -1109.5
$ws.Range("N7").Value = -2224
# Row 46
$ws.Range("H46").Value = 291141.56
$ws.Range("J46").Value = 6331.8335
$ws.Range("L46").Value = 6331.8335
$ws.Range("N46").Value = -6707.8335
# Row 93
$ws.Range("H93").Value = 1806.1428
$ws.Range("J93").Value = 397.5
$ws.Range("L93").Value = 397.5
$ws.Range("N93").Value = -2893.5
# Row 126
$ws.Range("H126").Value = 1377.2
$ws.Range("I126").Value = 1221.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3664.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1194.5
$ws.Range("N126").Value = -10940
# Row 136
$ws.Range("H136").Value = 4946
$ws.Range("I136").Value = 4999.6665
$ws.Range("K136").Value = 14998.9995
$ws.Range("M136").Value = -12448.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 6574.625
$ws.Range("I136").Value = 6574.625
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 19723.875
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -17173.875
$ws.Range("N136").ClearContents()

Write-Host "All edits applied"